# Auto-generated edits applying the market-price recompute diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3702.8
$ws.Range("I62").Value = 3378.5
$ws.Range("K62").Value = 3378.5
$ws.Range("M62").Value = -2754.5

$ws.Range("H65").Value = 3702.8
$ws.Range("I65").Value = 3378.5
$ws.Range("K65").Value = 16892.5
$ws.Range("M65").Value = -13772.5

$ws.Range("H70").Value = 3302.7058
$ws.Range("J70").Value = 4040.182
$ws.Range("L70").Value = 12120.546
$ws.Range("N70").Value = -12660.546

$ws.Range("H73").Value = 3302.7058
$ws.Range("J73").Value = 4040.182
$ws.Range("L73").Value = 12120.546
$ws.Range("N73").Value = -13992.546

$ws.Range("H75").Value = 69749.25
$ws.Range("J75").Value = 69749.25
$ws.Range("L75").Value = 69749.25
$ws.Range("N75").Value = -71621.25

$ws.Range("H78").Value = 69749.25
$ws.Range("J78").Value = 69749.25
$ws.Range("L78").Value = 209247.75
$ws.Range("N78").Value = -218607.75

$ws.Range("H137").Value = 4211.2896
$ws.Range("I137").Value = 4243.606
$ws.Range("K137").Value = 12730.818
$ws.Range("M137").Value = -10180.818

$ws.Range("H138").Value = 3482.8333
$ws.Range("I138").Value = 1936.7368
$ws.Range("K138").Value = 5810.2104
$ws.Range("M138").Value = -670.2103999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18523870
$ws.Range("I32").Value = 18523870
$ws.Range("K32").Value = 18523870
$ws.Range("M32").Value = -18523583

$ws.Range("H45").Value = 2093.52
$ws.Range("I45").Value = 917
$ws.Range("K45").Value = 917
$ws.Range("M45").Value = -540

$ws.Range("H61").Value = 2742.8125
$ws.Range("I61").Value = 2819
$ws.Range("K61").Value = 2819
$ws.Range("M61").Value = -2607

$ws.Range("H63").Value = 6379.75
$ws.Range("I63").Value = 3884.889
$ws.Range("J63").Value = 13864.333
$ws.Range("K63").Value = 3884.889
$ws.Range("L63").Value = 13864.333
$ws.Range("M63").Value = -3198.889
$ws.Range("N63").Value = -15236.333

$ws.Range("H66").Value = 6379.75
$ws.Range("I66").Value = 3884.889
$ws.Range("J66").Value = 13864.333
$ws.Range("K66").Value = 19424.445
$ws.Range("L66").Value = 69321.66500000001
$ws.Range("M66").Value = -15992.445
$ws.Range("N66").Value = -76185.66500000001

$ws.Range("H97").Value = 2634.9333
$ws.Range("I97").Value = 1310.3636
$ws.Range("K97").Value = 1310.3636
$ws.Range("M97").Value = -814.3635999999999

$ws.Range("H110").Value = 8006.143
$ws.Range("I110").Value = 8009.75
$ws.Range("J110").Value = 8001.3335
$ws.Range("K110").Value = 8009.75
$ws.Range("L110").Value = 8001.3335
$ws.Range("M110").Value = -5964.75
$ws.Range("N110").Value = -12091.3335

$ws.Range("H132").Value = 2835.0476
$ws.Range("I132").Value = 2807.111
$ws.Range("J132").Value = 3002.6667
$ws.Range("K132").Value = 8421.332999999999
$ws.Range("L132").Value = 9008.000100000001
$ws.Range("M132").Value = -5891.332999999999
$ws.Range("N132").Value = -14068.0001

$ws.Range("H136").Value = 2742.8125
$ws.Range("I136").Value = 2819
$ws.Range("K136").Value = 8457
$ws.Range("M136").Value = -5907

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 950.3333
$ws.Range("I64").Value = 785
$ws.Range("J64").Value = 1033
$ws.Range("K64").Value = 785
$ws.Range("L64").Value = 1033
$ws.Range("M64").Value = -560
$ws.Range("N64").Value = -1483

$ws.Range("H67").Value = 950.3333
$ws.Range("I67").Value = 785
$ws.Range("J67").Value = 1033
$ws.Range("K67").Value = 785
$ws.Range("L67").Value = 1033
$ws.Range("M67").Value = -5
$ws.Range("N67").Value = -2593

$ws.Range("H94").Value = 2232.5806
$ws.Range("I94").Value = 2768.9524
$ws.Range("J94").Value = 1106.2
$ws.Range("K94").Value = 2768.9524
$ws.Range("L94").Value = 1106.2
$ws.Range("M94").Value = -2317.9524
$ws.Range("N94").Value = -2008.2

$ws.Range("H99").Value = 100727.25
$ws.Range("I99").Value = 100727.25
$ws.Range("K99").Value = 100727.25
$ws.Range("M99").Value = -99229.25

$ws.Range("H134").Value = 2723.2
$ws.Range("I134").Value = 2521.3333
$ws.Range("K134").Value = 7563.999899999999
$ws.Range("M134").Value = -5028.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1899.7273
$ws.Range("I31").Value = 1857.2333
$ws.Range("K31").Value = 1857.2333
$ws.Range("M31").Value = -1562.2333

$ws.Range("H34").Value = 1899.7273
$ws.Range("I34").Value = 1857.2333
$ws.Range("K34").Value = 1857.2333
$ws.Range("M34").Value = -1655.2333

$ws.Range("H58").Value = 1312.409
$ws.Range("J58").Value = 3037.818
$ws.Range("L58").Value = 3037.818
$ws.Range("N58").Value = -3443.818

$ws.Range("H60").Value = 8499
$ws.Range("I60").Value = 8499
$ws.Range("K60").Value = 8499
$ws.Range("M60").Value = -7988

$ws.Range("H132").Value = 2132.4167
$ws.Range("I132").Value = 2132.4167
$ws.Range("K132").Value = 6397.250100000001
$ws.Range("M132").Value = -3867.250100000001

$ws.Range("H134").Value = 2528.0527
$ws.Range("I134").Value = 2183.5293
$ws.Range("K134").Value = 6550.5879
$ws.Range("M134").Value = -4015.5879

$ws.Range("H136").Value = 1312.409
$ws.Range("J136").Value = 3037.818
$ws.Range("L136").Value = 9113.454000000002
$ws.Range("N136").Value = -14213.454

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1541.5
$ws.Range("J113").Value = 1476.75
$ws.Range("L113").Value = 4430.25
$ws.Range("N113").Value = -8770.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4141.0938
$ws.Range("I122").Value = 3788.8948
$ws.Range("K122").Value = 11366.6844
$ws.Range("M122").Value = -8916.6844

$ws.Range("H132").Value = 2834.1
$ws.Range("I132").Value = 2668.25
$ws.Range("J132").Value = 3497.5
$ws.Range("K132").Value = 8004.75
$ws.Range("L132").Value = 10492.5
$ws.Range("M132").Value = -5474.75
$ws.Range("N132").Value = -15552.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3317.2354
$ws.Range("I68").Value = 3054.6667
$ws.Range("J68").Value = 3947.4
$ws.Range("K68").Value = 3054.6667
$ws.Range("L68").Value = 3947.4
$ws.Range("M68").Value = -2305.6667
$ws.Range("N68").Value = -5445.4

$ws.Range("H71").Value = 3317.2354
$ws.Range("I71").Value = 3054.6667
$ws.Range("J71").Value = 3947.4
$ws.Range("K71").Value = 15273.3335
$ws.Range("L71").Value = 19737
$ws.Range("M71").Value = -11529.3335
$ws.Range("N71").Value = -27225

$ws.Range("H74").Value = 70555.55499999999
$ws.Range("J74").Value = 76875
$ws.Range("L74").Value = 76875
$ws.Range("N74").Value = -78871

$ws.Range("H77").Value = 70555.55499999999
$ws.Range("J77").Value = 76875
$ws.Range("L77").Value = 230625
$ws.Range("N77").Value = -240609

$ws.Range("H122").Value = 5434.3706
$ws.Range("I122").Value = 2582.3
$ws.Range("J122").Value = 7112.0586
$ws.Range("K122").Value = 7746.900000000001
$ws.Range("L122").Value = 21336.1758
$ws.Range("M122").Value = -5296.900000000001
$ws.Range("N122").Value = -26236.1758

$ws.Range("H136").Value = 3486.8518
$ws.Range("I136").Value = 3428.6538
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 10285.9614
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -7735.9614
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H100").Value = 1054.6923
$ws.Range("I100").Value = 1128.2727
$ws.Range("K100").Value = 2256.5454
$ws.Range("M100").Value = -1715.5454

$ws.Range("H122").Value = 1973.1052
$ws.Range("I122").Value = 2016.3889
$ws.Range("K122").Value = 6049.1667
$ws.Range("M122").Value = -3599.1667

$ws.Range("H136").Value = 2106.3635
$ws.Range("I136").Value = 1096.4706
$ws.Range("K136").Value = 3289.4118
$ws.Range("M136").Value = -739.4118000000003
